# Mactaquac distribution workbook: add a "Trough" column to the Groups sheet
# (stock code report / group container history work), make Groups the active
# sheet/tab again, and update the helper comments that ride along the header
# row so they still line up with their columns.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Groups")

# --- 1. Remember the existing header comments for columns P..X (by their
#        current column letter) before anything moves. A column Insert shifts
#        cell values/styles automatically, but comments stay anchored to
#        their original cell, so they need to be re-homed by hand afterwards.
$sourceCols = @("P", "Q", "R", "S", "T", "U", "V", "X")
$savedText = @{}
foreach ($col in $sourceCols) {
    $cell = $ws.Range($col + "3")
    if ($cell.Comment -ne $null) {
        $savedText[$col] = $cell.Comment.Text()
    }
}

# --- 2. Insert the new "Trough" column at P, pushing Program..Exclude from
#        P:X out to Q:Y. ---
$ws.Columns("P").Insert()
$ws.Range("P3").Value = "Trough"

# give the new column roughly the same width as its "Year Collection" /
# "Program" neighbours (N:O) instead of the generic default width
$ws.Columns("P").ColumnWidth = $ws.Columns("O").ColumnWidth

# --- 3. Re-home the saved comments one column to the right. Columns Q..V
#        already have a comment object sitting on them (left over from the
#        pre-insert layout), so those just get their text swapped in place
#        -- processing right-to-left so a destination is never overwritten
#        before its old text has been read. W and Y don't have a comment
#        object yet, so those need a fresh one; X's now-stale comment (its
#        text moved to Y) is removed. ---
$null = $ws.Range("Y3").AddComment($savedText["X"])
$ws.Range("X3").Comment.Delete()

$null = $ws.Range("W3").AddComment($savedText["V"])

$ws.Range("V3").Comment.Text($savedText["U"])
$ws.Range("U3").Comment.Text($savedText["T"])
$ws.Range("T3").Comment.Text($savedText["S"])
$ws.Range("S3").Comment.Text($savedText["R"])
$ws.Range("R3").Comment.Text($savedText["Q"])
$ws.Range("Q3").Comment.Text($savedText["P"])

# --- 4. New comment for the new Trough column itself. ---
$null = $ws.Range("P3").AddComment("Enter trough names here if distributing from troughs.")

# --- 5. Groups becomes the active/selected sheet & tab again (it had lost
#        that to Individuals); selection sits on the new P3 cell. ---
$ws.Activate()
$null = $ws.Range("P3").Select()
